# Update "Jogos da Semana" FlashScore sheet:
#  - remove the last two scoreline-odds columns (Odd_CS_3-3_HT, Odd_CS_4-4_HT)
#  - replace the single data row with the new match/odds values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete columns BC:BD entirely (this also shifts the dimension to A1:BB2
# and removes both the "Odd_CS_3-3_HT"/"Odd_CS_4-4_HT" headers and their
# row-2 values automatically).
$ws.Range("BC1:BD2").EntireColumn.Delete()

# Row 2 - new match data / odds
$ws.Range("A2").Value  = "xjDW12Ug"
$ws.Range("C2").Value  = "09:00"
$ws.Range("D2").Value  = "THAILAND - THAI LEAGUE 1"
$ws.Range("E2").Value  = "Muang Thong Utd"
$ws.Range("F2").Value  = "Nakhon Pathom"
$ws.Range("G2").Value  = 1.3
$ws.Range("H2").Value  = 5.1
$ws.Range("I2").Value  = 8.75
$ws.Range("J2").Value  = 1.7
$ws.Range("K2").Value  = 2.6
$ws.Range("L2").Value  = 7.3
$ws.Range("N2").Value  = 9.5
$ws.Range("O2").Value  = 1.15
$ws.Range("P2").Value  = 4.75
$ws.Range("Q2").Value  = 1.47
$ws.Range("R2").Value  = 2.52
$ws.Range("S2").Value  = 1.27
$ws.Range("T2").Value  = 3.4
$ws.Range("U2").Value  = 1.78
$ws.Range("V2").Value  = 1.93
$ws.Range("W2").Value  = 9.25
$ws.Range("X2").Value  = 7.5
$ws.Range("Z2").Value  = 8.75
$ws.Range("AA2").Value = 10
$ws.Range("AB2").Value = 22
$ws.Range("AC2").Value = 9.5
$ws.Range("AD2").Value = 10.5
$ws.Range("AE2").Value = 19.5
$ws.Range("AF2").Value = 75
$ws.Range("AG2").Value = 500
$ws.Range("AH2").Value = 27
$ws.Range("AI2").Value = 65
$ws.Range("AJ2").Value = 27
$ws.Range("AK2").Value = 250
$ws.Range("AL2").Value = 100
$ws.Range("AM2").Value = 75
$ws.Range("AN2").Value = 3.3
$ws.Range("AO2").Value = 5.7
$ws.Range("AP2").Value = 14
$ws.Range("AQ2").Value = 14
$ws.Range("AR2").Value = 35
$ws.Range("AS2").Value = 150
$ws.Range("AT2").Value = 3.4
$ws.Range("AU2").Value = 8.25
$ws.Range("AV2").Value = 65
$ws.Range("AW2").Value = 9.75
$ws.Range("AX2").Value = 50
$ws.Range("AY2").Value = 45
$ws.Range("AZ2").Value = 350
$ws.Range("BA2").Value = 350
$ws.Range("BB2").Value = 500
